# Source "Patients" sheet selection changes to B2 (no longer the active tab)
$wb = $excel.ActiveWorkbook
$wsPatients = $wb.Worksheets.Item("Patients")
[void]$wsPatients.Range("B2").Select()

# "Providers" sheet selection changes to F57
$wsProviders = $wb.Worksheets.Item("Providers")
[void]$wsProviders.Range("F57").Select()

# Add a new worksheet "PatientsShifted" at the end, containing the same
# patient data as "Patients" but shifted one column right and five rows down
$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$ws.Name = "PatientsShifted"

# Pre-format the "Date of Birth" text cells as Text so the literal date-like
# strings are preserved instead of being auto-converted to date serials
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"

# Header row and data values (copy of "Patients"!A1:O5, shifted to B6:P10)
$ws.Range("B6").Value = "Patient ID"
$ws.Range("C6").Value = "Family Name"
$ws.Range("D6").Value = "Given Name"
$ws.Range("E6").Value = "Date of Birth"
$ws.Range("F6").Value = "Sex"
$ws.Range("G6").Value = "Contacts"
$ws.Range("H6").Value = "Admission Date"
$ws.Range("I6").Value = "Discharge Date"
$ws.Range("J6").Value = "Service Code"
$ws.Range("K6").Value = "Complete?"
$ws.Range("L6").Value = "BP Systolic"
$ws.Range("M6").Value = "BP Diastolic"
$ws.Range("N6").Value = "Weight (kg)"
$ws.Range("O6").Value = "Warfarin?"
$ws.Range("P6").Value = "Physician"
$ws.Range("B7").Value = 47234
$ws.Range("C7").Value = "Smith"
$ws.Range("D7").Value = "Jeff"
$ws.Range("E7").Value = "05/04/1968"
$ws.Range("F7").Value = "M"
$ws.Range("G7").Value = "ph: +15554441111"
$ws.Range("H7").Value = 45395
$ws.Range("I7").Value = 45427
$ws.Range("J7").Value = 73761001
$ws.Range("K7").Value = "y"
$ws.Range("L7").Value = 140
$ws.Range("M7").Value = 90
$ws.Range("N7").Value = 65
$ws.Range("O7").Value = "Y"
$ws.Range("P7").Value = "Werner von Braun"
$ws.Range("B8").Value = 689272
$ws.Range("C8").Value = "Brown"
$ws.Range("D8").Value = "Sue"
$ws.Range("E8").Value = "06/05/1972"
$ws.Range("F8").Value = "F"
$ws.Range("G8").Value = "em: sue@nowhere.com, ph: +6155443322"
$ws.Range("H8").Value = 45396
$ws.Range("I8").Value = 45462
$ws.Range("J8").Value = 26390003
$ws.Range("K8").Value = "y"
$ws.Range("L8").Value = 190
$ws.Range("M8").Value = 130
$ws.Range("O8").Value = "N"
$ws.Range("P8").Value = "Robert H Goddard"
$ws.Range("B9").Value = 2451
$ws.Range("C9").Value = "White"
$ws.Range("D9").Value = "David"
$ws.Range("E9").Value = "03/12/1955"
$ws.Range("F9").Value = "M"
$ws.Range("H9").Value = 45397
$ws.Range("I9").Value = 45474
$ws.Range("J9").Value = 726429001
$ws.Range("K9").Value = "y"
$ws.Range("L9").Value = 120
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 92
$ws.Range("P9").Value = "Sergei Korolev"
$ws.Range("B10").Value = 562
$ws.Range("C10").Value = "Green"
$ws.Range("D10").Value = "Kath"
$ws.Range("E10").Value = "5/23/1991"
$ws.Range("F10").Value = "F"
$ws.Range("G10").Value = "fax: +15553332222"
$ws.Range("H10").Value = 45399
$ws.Range("J10").Value = 39633000
$ws.Range("K10").Value = "n"
$ws.Range("L10").Value = 123
$ws.Range("M10").Value = 76
$ws.Range("N10").Value = 36
$ws.Range("O10").Value = "?"
$ws.Range("P10").Value = "Konstantin Tisolkovsky"

# Apply number formats matching the source columns
# Admission Date column
$ws.Range("H7").NumberFormat = "m/d/yy"
$ws.Range("H8").NumberFormat = "m/d/yy"
$ws.Range("H9").NumberFormat = "m/d/yy"
$ws.Range("H10").NumberFormat = "m/d/yy"
# Discharge Date column
$ws.Range("I7").NumberFormat = "d-mmm-yy"
$ws.Range("I8").NumberFormat = "d-mmm-yy"
$ws.Range("I9").NumberFormat = "d-mmm-yy"
$ws.Range("I10").NumberFormat = "d-mmm-yy"
# Date of Birth column
$ws.Range("E7").NumberFormat = "mm/dd/yy;@"
$ws.Range("E8").NumberFormat = "mm/dd/yy;@"
$ws.Range("E9").NumberFormat = "mm/dd/yy;@"
$ws.Range("E10").NumberFormat = "mm/dd/yy;@"
# Date of Birth header
$ws.Range("E6").NumberFormat = "@"

# Empty cell that still carries the Discharge Date number format
$ws.Range("I10").NumberFormat = "d-mmm-yy"

# Select G7 on the new sheet, making it the active tab (matches the recorded edit)
[void]$ws.Range("G7").Select()
